$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "app"
$ws.Range("B3").Value = "cameron"
$ws.Range("C3").Value = "testing app feedback"
$ws.Range("D3").Value = "2025-10-01 16:39:21"
